$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Versions")

# Fill in newly-added version values for row 18
$ws.Range("D18").Value = "6.0.2"
$ws.Range("F18").Value = "10.1.0"
$ws.Range("H18").Value = "12.1.0"

# Update the active selection on the sheet
$ws.Range("H18").Select()
